$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a literal text value (so date-looking strings like
# "10/30/2022" are stored as text, same as the existing shared-string
# dates in this sheet) instead of being auto-converted into a date serial.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Push back the release dates for the two Estimation Statistics sessions
# by one week (10/23/2022 -> 10/30/2022)...
Set-TextValue "C9" "10/30/2022"
Set-TextValue "C10" "10/30/2022"

# ...and the Fundamental Forecasting Models session moves out a week too
# (10/30/2022 -> 11/06/2022).
Set-TextValue "C12" "11/06/2022"

# Update the saved cursor/selection position on the sheet.
$ws.Range("D17").Select()
